$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 319.18182
$ws.Range("I2").Value = 321.1
$ws.Range("K2").Value = 321.1
$ws.Range("M2").Value = -208.1
$ws.Range("H11").Value = 27
$ws.Range("I11").Value = 27
$ws.Range("K11").Value = 27
$ws.Range("M11").Value = 113
$ws.Range("H70").Value = 2834.8235
$ws.Range("J70").Value = 3090.3635
$ws.Range("L70").Value = 9271.0905
$ws.Range("N70").Value = -9811.0905
$ws.Range("H73").Value = 2834.8235
$ws.Range("J73").Value = 3090.3635
$ws.Range("L73").Value = 9271.0905
$ws.Range("N73").Value = -11143.0905
$ws.Range("H80").Value = 3584.1
$ws.Range("I80").Value = 2120.25
$ws.Range("K80").Value = 6360.75
$ws.Range("M80").Value = -5362.75
$ws.Range("H83").Value = 3584.1
$ws.Range("I83").Value = 2120.25
$ws.Range("K83").Value = 19082.25
$ws.Range("M83").Value = -14090.25
$ws.Range("H112").Value = 1661.091
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1661.091
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4983.272999999999
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -7199.272999999999
$ws.Range("H114").Value = 52999
$ws.Range("J114").Value = 52999
$ws.Range("L114").Value = 52999
$ws.Range("N114").Value = -61677
$ws.Range("H125").Value = 803.6
$ws.Range("I125").Value = 499
$ws.Range("K125").Value = 4491
$ws.Range("M125").Value = -2031
$ws.Range("H127").Value = 1396.25
$ws.Range("I127").Value = 1454.2858
$ws.Range("K127").Value = 4362.857400000001
$ws.Range("M127").Value = 597.1425999999992
$ws.Range("H135").Value = 722.875
$ws.Range("I135").Value = 822
$ws.Range("K135").Value = 7398
$ws.Range("M135").Value = -4863

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1571.1333
$ws.Range("I45").Value = 1480.75
$ws.Range("J45").Value = 1932.6666
$ws.Range("K45").Value = 1480.75
$ws.Range("L45").Value = 1932.6666
$ws.Range("M45").Value = -1103.75
$ws.Range("N45").Value = -2686.6666
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H110").Value = 5246.1177
$ws.Range("I110").Value = 5378.4546
$ws.Range("J110").Value = 5003.5
$ws.Range("K110").Value = 5378.4546
$ws.Range("L110").Value = 5003.5
$ws.Range("M110").Value = -3333.4546
$ws.Range("N110").Value = -9093.5
$ws.Range("H132").Value = 9618959
$ws.Range("I132").Value = 2398
$ws.Range("J132").Value = 55564750
$ws.Range("K132").Value = 7194
$ws.Range("L132").Value = 166694250
$ws.Range("M132").Value = -4664
$ws.Range("N132").Value = -166699310

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 565.7692
$ws.Range("I22").Value = 541
$ws.Range("K22").Value = 541
$ws.Range("M22").Value = -368
$ws.Range("H99").Value = 1827.7142
$ws.Range("I99").Value = 1821.5555
$ws.Range("K99").Value = 1821.5555
$ws.Range("M99").Value = -323.5554999999999
$ws.Range("H105").Value = 2661.524
$ws.Range("I105").Value = 2637.8572
$ws.Range("J105").Value = 2708.8572
$ws.Range("K105").Value = 2637.8572
$ws.Range("L105").Value = 2708.8572
$ws.Range("M105").Value = -890.8571999999999
$ws.Range("N105").Value = -6202.8572

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 17857.428
$ws.Range("I62").Value = 10999.333
$ws.Range("K62").Value = 10999.333
$ws.Range("M62").Value = -10375.333
$ws.Range("H65").Value = 17857.428
$ws.Range("I65").Value = 10999.333
$ws.Range("K65").Value = 54996.665
$ws.Range("M65").Value = -51876.665

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 496.83334
$ws.Range("I5").Value = 496.83334
$ws.Range("K5").Value = 1490.50002
$ws.Range("M5").Value = -1378.50002
$ws.Range("H11").Value = 8121.857
$ws.Range("I11").Value = 9169.799999999999
$ws.Range("J11").Value = 5502
$ws.Range("K11").Value = 27509.4
$ws.Range("L11").Value = 16506
$ws.Range("M11").Value = -27369.4
$ws.Range("N11").Value = -16786
$ws.Range("H23").Value = 2225
$ws.Range("I23").Value = 2835.4
$ws.Range("K23").Value = 8506.200000000001
$ws.Range("M23").Value = -8271.200000000001
$ws.Range("H80").Value = 4461.4546
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 4457.6
$ws.Range("K80").Value = 13500
$ws.Range("L80").Value = 13372.8
$ws.Range("M80").Value = -12564
$ws.Range("N80").Value = -15244.8
$ws.Range("H81").Value = 6887.3335
$ws.Range("I81").Value = 6284.143
$ws.Range("K81").Value = 18852.429
$ws.Range("M81").Value = -17729.429
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H83").Value = 4461.4546
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 4457.6
$ws.Range("K83").Value = 40500
$ws.Range("L83").Value = 40118.4
$ws.Range("M83").Value = -35820
$ws.Range("N83").Value = -49478.4
$ws.Range("H84").Value = 6887.3335
$ws.Range("I84").Value = 6284.143
$ws.Range("K84").Value = 56557.287
$ws.Range("M84").Value = -50941.287
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H88").Value = 2699
$ws.Range("I88").Value = 2499
$ws.Range("K88").Value = 7497
$ws.Range("M88").Value = -7069
$ws.Range("H91").Value = 2699
$ws.Range("I91").Value = 2499
$ws.Range("K91").Value = 7497
$ws.Range("M91").Value = -6015
$ws.Range("H104").Value = 3500
$ws.Range("J104").Value = 3500
$ws.Range("L104").Value = 10500
$ws.Range("N104").Value = -15742
$ws.Range("H107").Value = 419
$ws.Range("J107").Value = 448.75
$ws.Range("L107").Value = 1346.25
$ws.Range("N107").Value = -5186.25
$ws.Range("H135").Value = 496.83334
$ws.Range("I135").Value = 496.83334
$ws.Range("K135").Value = 4471.50006
$ws.Range("M135").Value = -1936.50006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 932.4783
$ws.Range("I107").Value = 984.9
$ws.Range("J107").Value = 583
$ws.Range("K107").Value = 984.9
$ws.Range("L107").Value = 583
$ws.Range("M107").Value = 935.1
$ws.Range("N107").Value = -4423

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2093.3
$ws.Range("I16").Value = 1366.625
$ws.Range("K16").Value = 1366.625
$ws.Range("M16").Value = -1196.625
$ws.Range("H22").Value = 3462.818
$ws.Range("I22").Value = 2022.1666
$ws.Range("J22").Value = 5191.6
$ws.Range("K22").Value = 2022.1666
$ws.Range("L22").Value = 5191.6
$ws.Range("M22").Value = -1727.1666
$ws.Range("N22").Value = -5781.6
$ws.Range("H27").Value = 3462.818
$ws.Range("I27").Value = 2022.1666
$ws.Range("J27").Value = 5191.6
$ws.Range("K27").Value = 2022.1666
$ws.Range("L27").Value = 5191.6
$ws.Range("M27").Value = -1915.1666
$ws.Range("N27").Value = -5405.6
$ws.Range("H29").Value = 15995
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 15995
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 15995
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -16585
$ws.Range("H61").Value = 2972
$ws.Range("I61").Value = 2967.75
$ws.Range("K61").Value = 2967.75
$ws.Range("M61").Value = -2765.75
$ws.Range("H113").Value = 2972
$ws.Range("I113").Value = 2967.75
$ws.Range("K113").Value = 2967.75
$ws.Range("M113").Value = -797.75
$ws.Range("H136").Value = 30306858
$ws.Range("I136").Value = 2745.3333
$ws.Range("K136").Value = 8235.999899999999
$ws.Range("M136").Value = -5685.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 937.1429000000001
$ws.Range("I136").Value = 937.1429000000001
$ws.Range("K136").Value = 2811.4287
$ws.Range("M136").Value = -261.4287000000004
